$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "jo"
$ws.Range("B2").Value = 3

$ws.Range("A3").Value = "ji"
$ws.Range("B3").Value = 4

$ws.Range("A4").Value = "jp"
$ws.Range("B4").Value = 2
